# Generate Report for Handback
#
# Fills in the "Latest Target File" (F) and "Latest Handback File" (G)
# hyperlink columns for the zh-cn and de-de handoff/handback sheets, marks
# the two tracked docs as handed back (in sync with en-US), and stamps the
# handback datetime for each locale.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns (zh-cn / de-de) for both tracked docs ---
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Range("C3").Value = $statusHandedBack

# Row 2 (5c821a51-fe42-436c-819f-a1c9eeb68368 doc): target + handback files,
# pointing at the same source the handoff used (in sync -> same file).
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/97d6ff8454d89a64d17bb4a9c26b20aed888c720/e2e/5c821a51-fe42-436c-819f-a1c9eeb68368.md",
    "",
    "",
    "5c821a51-fe42-436c-819f-a1c9eeb68368.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5064731d7a8170962468770d4bf592c9112820a8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5c821a51-fe42-436c-819f-a1c9eeb68368.733048542cdf4747ff5fec3347bd913d8d690fba.zh-cn.xlf",
    "",
    "",
    "5c821a51-fe42-436c-819f-a1c9eeb68368.733048542cdf4747ff5fec3347bd913d8d690fba.zh-cn.xlf"
)

# Row 3 (96ea4119-27f9-41b6-9d33-b6b848f64680 doc)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/97d6ff8454d89a64d17bb4a9c26b20aed888c720/e2e/96ea4119-27f9-41b6-9d33-b6b848f64680.md",
    "",
    "",
    "96ea4119-27f9-41b6-9d33-b6b848f64680.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5064731d7a8170962468770d4bf592c9112820a8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/96ea4119-27f9-41b6-9d33-b6b848f64680.f95422c26a1d8c8722209a45d6b3c2469c3773cf.zh-cn.xlf",
    "",
    "",
    "96ea4119-27f9-41b6-9d33-b6b848f64680.f95422c26a1d8c8722209a45d6b3c2469c3773cf.zh-cn.xlf"
)

# zh-cn handback happened at 2016-03-21 12:40:15
$wsZhCn.Range("H2").Value = "2016-03-21 12:40:15"
$wsZhCn.Range("H3").Value = "2016-03-21 12:40:15"

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Range("C3").Value = $statusHandedBack

# Row 2 (5c821a51-fe42-436c-819f-a1c9eeb68368 doc)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/97d6ff8454d89a64d17bb4a9c26b20aed888c720/e2e/5c821a51-fe42-436c-819f-a1c9eeb68368.md",
    "",
    "",
    "5c821a51-fe42-436c-819f-a1c9eeb68368.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/045c75290d96f167b4400939621e025a90bb8717/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5c821a51-fe42-436c-819f-a1c9eeb68368.733048542cdf4747ff5fec3347bd913d8d690fba.de-de.xlf",
    "",
    "",
    "5c821a51-fe42-436c-819f-a1c9eeb68368.733048542cdf4747ff5fec3347bd913d8d690fba.de-de.xlf"
)

# Row 3 (96ea4119-27f9-41b6-9d33-b6b848f64680 doc)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/97d6ff8454d89a64d17bb4a9c26b20aed888c720/e2e/96ea4119-27f9-41b6-9d33-b6b848f64680.md",
    "",
    "",
    "96ea4119-27f9-41b6-9d33-b6b848f64680.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/045c75290d96f167b4400939621e025a90bb8717/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/96ea4119-27f9-41b6-9d33-b6b848f64680.f95422c26a1d8c8722209a45d6b3c2469c3773cf.de-de.xlf",
    "",
    "",
    "96ea4119-27f9-41b6-9d33-b6b848f64680.f95422c26a1d8c8722209a45d6b3c2469c3773cf.de-de.xlf"
)

# de-de handback happened at 2016-03-21 12:40:23
$wsDeDe.Range("H2").Value = "2016-03-21 12:40:23"
$wsDeDe.Range("H3").Value = "2016-03-21 12:40:23"
